# Remove D9, D10, R13, and R14 from the BOM
#
# LEDs D9 and D10 put too much load on the level shifter U4. Remove the
# LEDs (and the now-stale Rev-B changelog entries covering resistors
# R13/R14 and friends) from the Bill of Materials sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")

# --- D2/D9/D10 LED row: drop D9 and D10, leaving only D2 ---------------
$ws.Range("A13").Value = "D2"
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = "PWR"

# --- R2/R3/R4/R13/R14 resistor row: drop R13 and R14 --------------------
$ws.Range("A25").Value = "R2 R3 R4"
$ws.Range("B25").Value = 3

# --- M2.5 screw quantity bump (4 -> 8) ----------------------------------
$ws.Range("B38").Value = 8

# --- Trim the "Changes from Rev B" changelog block ----------------------
# Keep the "Note" / "Select components..." / "Seeed OPL components..."
# rows, but drop the heading and all of the stale per-revision bullets.
$ws.Range("A41").ClearContents()
$ws.Range("A42").ClearContents()
$ws.Range("A43:H49").EntireRow.Delete()
